$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

# Student B (row 3): Bonus (H3) corrected from 20 to 0, Total (C3) recalculated
$ws.Range("H3").Value = 0
$ws.Range("C3").Value = "27 (87.1%)"

# Student D (row 5): Task 2 (E5) corrected from 0 to 5, Total (C5) and Mark (B5) recalculated
$ws.Range("E5").Value = 5
$ws.Range("C5").Value = "14 (45.2%)"

# Set Mark (B5) as text "4" without leaving a stray number format / style behind
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "4"
$ws.Range("B5").ClearFormats()
